$wb = $excel.ActiveWorkbook

$wsItems = $wb.Worksheets.Item("Items")
$wsBestiary = $wb.Worksheets.Item("Bestiary")
$wsMoves = $wb.Worksheets.Item("Sample_Custom_Moves")

# --- Add two new Bestiary entries (rows 18 and 19) ---
$wsBestiary.Range("A18").Value = "Waterbone"
$wsBestiary.Range("B18").Value = "A skeleton inside a bubble of water"
$wsBestiary.Range("C18").Value = 9
$wsBestiary.Range("D18").Value = 1
$wsBestiary.Range("E18").Value = "Holy, Lightning"
$wsBestiary.Range("F18").Value = "Bubble Attack: 1d10 Water damage`nBone Attack: 1d12 Dark damage."

$wsBestiary.Range("A19").Value = "Electric Eel"
$wsBestiary.Range("B19").Value = "An underwater creature that shoots electricity"
$wsBestiary.Range("C19").Value = 9
$wsBestiary.Range("D19").Value = 0
$wsBestiary.Range("F19").Value = "Bite: 1d10 damage`nShock: 1d8 Lightning damage to all enemies.  Causes Paralysis."

$wsBestiary.Rows.Item(18).RowHeight = 45
$wsBestiary.Rows.Item(19).RowHeight = 45

# --- Update sheet view selections / active sheet ---
# Final active sheet should be Bestiary (activeTab=1), with topLeftCell A9, selection F20
$wsItems.Range("K25").Select() | Out-Null
$wsMoves.Range("A38").Select() | Out-Null
$wsBestiary.Range("F20").Select() | Out-Null
